$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly-scoreboard rows appended by Steven (Brave Leopard), week 7.
# Row 287: Walk on 7/26/2024
# Row 288: Run on 7/27/2024
# Row 289: Workout on 7/27/2024
$newRows = @(
    @{ Row = 287; Participant = "Steven"; Date = 45499; Type = "Walk";    Duration = 21; Distance = 1.1;  Elevation = 30; Z1 = 21; Z2 = 0; Z3 = 0; Z4 = 0; Z5 = 0; Level = "Brave Leopard"; Week = 7 },
    @{ Row = 288; Participant = "Steven"; Date = 45500; Type = "Run";     Duration = 13; Distance = 1.38; Elevation = 66; Z1 = 1;  Z2 = 1; Z3 = 8; Z4 = 3; Z5 = 0; Level = "Brave Leopard"; Week = 7 },
    @{ Row = 289; Participant = "Steven"; Date = 45500; Type = "Workout"; Duration = 10; Distance = 0;    Elevation = 0;  Z1 = 5;  Z2 = 5; Z3 = 0; Z4 = 0; Z5 = 0; Level = "Brave Leopard"; Week = 7 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value2 = $r.Participant

    # Copy the date cell's number format (s="1", m/d/yyyy) from the row above,
    # then overwrite just the value.
    $ws.Cells.Item($row - 1, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value2 = $r.Date

    $ws.Cells.Item($row, 3).Value2 = $r.Type
    $ws.Cells.Item($row, 4).Value2 = $r.Duration
    $ws.Cells.Item($row, 5).Value2 = $r.Distance
    $ws.Cells.Item($row, 6).Value2 = $r.Elevation
    $ws.Cells.Item($row, 7).Value2 = $r.Z1
    $ws.Cells.Item($row, 8).Value2 = $r.Z2
    $ws.Cells.Item($row, 9).Value2 = $r.Z3
    $ws.Cells.Item($row, 10).Value2 = $r.Z4
    $ws.Cells.Item($row, 11).Value2 = $r.Z5
    $ws.Cells.Item($row, 12).Value2 = $r.Level
    $ws.Cells.Item($row, 13).Value2 = $r.Week
}

$excel.CutCopyMode = $false

# Reflect the author's final cursor position: next blank row under the new data.
$ws.Range("A290").Select() | Out-Null
